$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 166 - this shifts the former rows 166-200 down to 167-201
$ws.Rows(166).Insert()

# Populate the newly inserted row 166 with the new weekly record.
# (Columns A, B, C, E, F, G, H, I, J, Q, R, T repeat the constant values
#  used throughout this block of rows; D, K, L, M, N, O, P, S are the new data.)
$ws.Range("A166").Value = 1
$ws.Range("B166").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C166").Value = "Arica y Parinacota"
$ws.Range("D166").Value = 44951
$ws.Range("E166").Value = 15
$ws.Range("F166").Value = "Fruta"
$ws.Range("G166").Value = 100108
$ws.Range("H166").Value = "Tropicales y subtropicales"
$ws.Range("I166").Value = 100108002
$ws.Range("J166").Value = "Mango"
$ws.Range("K166").Value = "Piqueño"
$ws.Range("L166").Value = "Primera"
$ws.Range("M166").Value = 456
$ws.Range("N166").Value = 4000
$ws.Range("O166").Value = 4500
$ws.Range("P166").Value = 4250
$ws.Range("Q166").Value = "$/bandeja 4 kilos"
$ws.Range("R166").Value = "Perú"
$ws.Range("S166").Value = 1062
$ws.Range("T166").Value = 4
